# Word COM-interop edit script.
# Source: answers-of-two-digit_number_divided_by_one-digit_number.docx
# Updates the date heading and the 25 division-problem answers in the table.

$d = $word.ActiveDocument

# --- 1. Update the date heading paragraph ("2025-03-24 Monday" -> "2025-03-25 Tuesday") ---
# The heading text is unique in the document, so Find/Replace is unambiguous
# and keeps the existing run formatting (Arial, sz 30) untouched.
$d.Content.Find.Execute("2025-03-24 Monday", $false, $false, $false, $false, $false, $true, 1, $false, "2025-03-25 Tuesday", 2) | Out-Null

# --- 2. Update the 25 answers in the single table ---
# Several old/new values collide across cells (e.g. the old text of one cell
# equals the new text of another), so a blanket Find/Replace over the whole
# document could clobber the wrong cell once the first replacement lands.
# Addressing every cell explicitly by (row, column) sidesteps that entirely.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "90÷7=12, 6"  # was "53÷9=5, 8"
$t.Cell(1, 2).Range.Text = "25÷8=3, 1"  # was "54÷7=7, 5"
$t.Cell(1, 3).Range.Text = "97÷9=10, 7"  # was "62÷4=15, 2"
$t.Cell(1, 4).Range.Text = "46÷4=11, 2"  # was "63÷6=10, 3"
$t.Cell(1, 5).Range.Text = "28÷4=7, 0"  # was "17÷6=2, 5"
$t.Cell(5, 1).Range.Text = "15÷3=5, 0"  # was "12÷8=1, 4"
$t.Cell(5, 2).Range.Text = "33÷9=3, 6"  # was "41÷2=20, 1"
$t.Cell(5, 3).Range.Text = "24÷9=2, 6"  # was "29÷8=3, 5"
$t.Cell(5, 4).Range.Text = "40÷6=6, 4"  # was "64÷8=8, 0"
$t.Cell(5, 5).Range.Text = "45÷3=15, 0"  # was "96÷3=32, 0"
$t.Cell(9, 1).Range.Text = "15÷6=2, 3"  # was "72÷7=10, 2"
$t.Cell(9, 2).Range.Text = "35÷5=7, 0"  # was "90÷7=12, 6"
$t.Cell(9, 3).Range.Text = "92÷9=10, 2"  # was "69÷4=17, 1"
$t.Cell(9, 4).Range.Text = "33÷8=4, 1"  # was "91÷4=22, 3"
$t.Cell(9, 5).Range.Text = "11÷6=1, 5"  # was "20÷5=4, 0"
$t.Cell(13, 1).Range.Text = "64÷2=32, 0"  # was "22÷2=11, 0"
$t.Cell(13, 2).Range.Text = "65÷4=16, 1"  # was "34÷8=4, 2"
$t.Cell(13, 3).Range.Text = "73÷6=12, 1"  # was "35÷2=17, 1"
$t.Cell(13, 4).Range.Text = "39÷3=13, 0"  # was "27÷2=13, 1"
$t.Cell(13, 5).Range.Text = "90÷4=22, 2"  # was "18÷6=3, 0"
$t.Cell(17, 1).Range.Text = "22÷4=5, 2"  # was "54÷3=18, 0"
$t.Cell(17, 2).Range.Text = "77÷7=11, 0"  # was "64÷4=16, 0"
$t.Cell(17, 3).Range.Text = "15÷7=2, 1"  # was "58÷2=29, 0"
$t.Cell(17, 4).Range.Text = "20÷9=2, 2"  # was "25÷7=3, 4"
$t.Cell(17, 5).Range.Text = "44÷2=22, 0"  # was "20÷8=2, 4"
